# Apply the Dec-11-2023 cryptos list refresh (prices + 1h volume deltas).
# Row 10/11 also swap places (Avalanche <-> Dogecoin) to mirror the new ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '40.666.06'
$ws.Range("E2").Value = '  -7.35%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.185.64'
$ws.Range("E3").Value = '  -7.37%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.08%  '

# Row 5: BNB
$ws.Range("D5").Value = '''241.58'
$ws.Range("E5").Value = '  +0.56%  '

# Row 6: XRP
$ws.Range("E6").Value = '  -7.91%  '

# Row 7: Solana
$ws.Range("D7").Value = '''68.27'
$ws.Range("E7").Value = '  -7.76%  '

# Row 8: USDC
$ws.Range("E8").Value = '  +0.19%  '

# Row 9: Cardano
$ws.Range("E9").Value = '  -12.37%  '

# Row 10: Dogecoin (was Avalanche)
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '''0.0940'
$ws.Range("E10").Value = '  -8.14%  '

# Row 11: Avalanche (was Dogecoin)
$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").Value = '''36.17'
$ws.Range("E11").Value = '  -2.46%  '

# Row 12: OKB
$ws.Range("D12").Value = '''57.52'
$ws.Range("E12").Value = '  -5.40%  '

# Row 13: TRON
$ws.Range("D13").Value = '''0.104'
$ws.Range("E13").Value = '  -4.53%  '

# Row 14: Polkadot
$ws.Range("D14").Value = '''6.54'
$ws.Range("E14").Value = '  -10.15%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '2.512.53'
$ws.Range("E15").Value = '  -7.38%  '

# Row 16: Chainlink
$ws.Range("E16").Value = '  -10.85%  '

# Row 17: Polygon
$ws.Range("D17").Value = '''0.826'

# Row 18: WrappedEther
$ws.Range("D18").Value = '2.187.84'
$ws.Range("E18").Value = '  -7.40%  '

# Row 19: WrappedBTC
$ws.Range("D19").Value = '40.681.18'
$ws.Range("E19").Value = '  -7.19%  '

# Row 20: ShibaInu
$ws.Range("D20").Value = '0.0₃0937'
$ws.Range("E20").Value = '  -9.27%  '

# Row 21: Litecoin
$ws.Range("D21").Value = '''72.12'
$ws.Range("E21").Value = '  -7.48%  '

# Row 22: Uniswap
$ws.Range("E22").Value = '  -8.24%  '

# Row 23: BitcoinCash
$ws.Range("D23").Value = '''228.58'
$ws.Range("E23").Value = '  -9.92%  '

# Row 24: ImmutableX
$ws.Range("E24").Value = '  +6.54%  '

# Row 25: Dai
$ws.Range("D25").Value = '''1.00'
$ws.Range("E25").Value = '  +0.07%  '

# Row 26: WEMIXToken
$ws.Range("E26").Value = '  -5.30%  '

# Row 27: PancakeSwap
$ws.Range("E27").Value = '  -4.41%  '

# Row 28: Toncoin
$ws.Range("E28").Value = '  -5.27%  '

# Row 29: Cosmos
$ws.Range("D29").Value = '''9.61'
$ws.Range("E29").Value = '  -8.59%  '

# Row 30: Monero
$ws.Range("D30").Value = '''168.35'
$ws.Range("E30").Value = '  -4.13%  '

# Row 31: EthereumClassic
$ws.Range("D31").Value = '''20.10'
$ws.Range("E31").Value = '  -10.07%  '

# Row 32: Kaspa
$ws.Range("E32").Value = '  -10.17%  '

# Row 33: Stellar
$ws.Range("E33").Value = '  -8.43%  '

# Row 34: Hedera
$ws.Range("D34").Value = '''0.0695'
$ws.Range("E34").Value = '  -7.65%  '

# Row 35: InternetComputer(DFINITY)
$ws.Range("E35").Value = '  -5.45%  '

# Row 36: Filecoin
$ws.Range("E36").Value = '  -10.66%  '

# Row 37: RenderToken
$ws.Range("D37").Value = '''3.81'
$ws.Range("E37").Value = '  -0.06%  '

# Row 38: InjectiveProtocol
$ws.Range("D38").Value = '''23.30'
$ws.Range("E38").Value = '  +14.42%  '

# Row 39: LidoDAOToken
$ws.Range("E39").Value = '  -6.71%  '

# Row 40: VeChain
$ws.Range("D40").Value = '''0.0267'
$ws.Range("E40").Value = '  -4.31%  '

# Row 41: THORChain
$ws.Range("E41").Value = '  -12.90%  '

# Row 42: MultiversX
$ws.Range("E42").Value = '  -4.66%  '

# Row 43: FTXToken
$ws.Range("D43").Value = '''4.80'
$ws.Range("E43").Value = '  -11.50%  '

# Row 44: FraxShare
$ws.Range("D44").Value = '''8.54'
$ws.Range("E44").Value = '  -5.80%  '

# Row 45: Algorand
$ws.Range("E45").Value = '  -7.19%  '

# Row 46: BinanceUSD
$ws.Range("E46").Value = '  +0.17%  '

# Row 47: Cronos
$ws.Range("D47").Value = '''0.0977'
$ws.Range("E47").Value = '  -8.29%  '

# Row 48: SynthetixNetwork
$ws.Range("D48").Value = '''4.43'
$ws.Range("E48").Value = '  +0.68%  '

# Row 49: Celestia
$ws.Range("D49").Value = '''10.20'
$ws.Range("E49").Value = '  +5.28%  '

# Row 50: TrustWalletToken
$ws.Range("E50").Value = '  -7.14%  '

# Row 51: ARBITRUM
$ws.Range("E51").Value = '  -6.66%  '
